$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) label renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) recalculated values ---
$ws.Range("C2").Value = 2934.187009790061
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1460.056109840828
$ws.Range("C5").Value = 5191.140356354663
$ws.Range("C6").Value = 4729.735976516416
$ws.Range("C7").Value = 2812.435974421079
$ws.Range("C8").Value = 401.8350013668368
$ws.Range("C9").Value = 951.6879611168786
$ws.Range("C10").Value = 2983.242707849043
$ws.Range("C11").Value = 2898.942214704482
$ws.Range("C12").Value = 665.6274194933962
$ws.Range("C13").Value = 1503.870423231357
$ws.Range("C14").Value = 5555.389721901988
$ws.Range("C15").Value = 5082.354756663512
$ws.Range("C16").Value = 2828.483778716848
$ws.Range("C17").Value = 1132.548400540401
$ws.Range("C18").Value = 417.6031683854853
$ws.Range("C19").Value = 982.980837581714
$ws.Range("C20").Value = 864.5379000312432
$ws.Range("C21").Value = 3083.80337578809
$ws.Range("C22").Value = 2965.153206179127
$ws.Range("C23").Value = 691.8942672110555
$ws.Range("C24").Value = 1577.487171555845
$ws.Range("C25").Value = 5660.517066940175
$ws.Range("C26").Value = 4921.848409120176
$ws.Range("C27").Value = 1657.651524528445
$ws.Range("C28").Value = 2999.422762626143
$ws.Range("C29").Value = 3156.723844635973
$ws.Range("C30").Value = 5122.180090208862
$ws.Range("C31").Value = 5642.578115155247
$ws.Range("C32").Value = 2860.874335573629
$ws.Range("C33").Value = 1000.829216794104
$ws.Range("C34").Value = 0
$ws.Range("C35").Value = 5745.422744292303
$ws.Range("C36").Value = 701.4459636783288
$ws.Range("C37").Value = 869.6014949562591
$ws.Range("C39").Value = 1716.389195271215
$ws.Range("C40").Value = 3056.152683606517
$ws.Range("C41").Value = 3212.740625904757
$ws.Range("C42").Value = 5295.682695961288
$ws.Range("C43").Value = 5919.20956823756
$ws.Range("C44").Value = 2887.250212489506
$ws.Range("C45").Value = 1032.277326842402
$ws.Range("C46").Value = 0
$ws.Range("C47").Value = 5955.175904294275
$ws.Range("C48").Value = 720.7128711178943
$ws.Range("C49").Value = 872.1235974568563
$ws.Range("C51").Value = 3008.669179463094
$ws.Range("C52").Value = 3252.634165082374
$ws.Range("C53").Value = 3137.260298393558
$ws.Range("C54").Value = 730.3063521039821
$ws.Range("C55").Value = 1060.095015975378
$ws.Range("C56").Value = 707.8672001573369
$ws.Range("C57").Value = 711.3043470146426
$ws.Range("C58").Value = 1775.027517189621
$ws.Range("C59").Value = 5996.49696468919
$ws.Range("C61").Value = 6301.696269820412
$ws.Range("C62").Value = 6103.744960203087
$ws.Range("C63").Value = 0
$ws.Range("C65").Value = 3012.536723186288
$ws.Range("C66").Value = 3314.741082534716
$ws.Range("C67").Value = 3210.869677115934
$ws.Range("C68").Value = 729.1196658666737
$ws.Range("C69").Value = 1093.134170274031
$ws.Range("C70").Value = 729.7808175407341
$ws.Range("C71").Value = 731.9993357350996
$ws.Range("C72").Value = 1836.014008604312
$ws.Range("C73").Value = 6114.227214287786
$ws.Range("C75").Value = 6661.86504232374
$ws.Range("C76").Value = 6249.151036691844
$ws.Range("C77").Value = 0
$ws.Range("C79").Value = 2854.757682901436
$ws.Range("C80").Value = 5176.058803160127
$ws.Range("C81").Value = 3382.563653843273
$ws.Range("C82").Value = 3242.636921959078
$ws.Range("C83").Value = 729.8559996981501
$ws.Range("C84").Value = 1129.713195979213
$ws.Range("C85").Value = 749.2194349876407
$ws.Range("C86").Value = 1895.214690888655
$ws.Range("C87").Value = 6262.368904654469
$ws.Range("C88").Value = 0

# --- Column AL flag flips (0 -> 1) ---
$ws.Range("AL5").Value = 1
$ws.Range("AL12").Value = 1
$ws.Range("AL14").Value = 1
$ws.Range("AL23").Value = 1
$ws.Range("AL25").Value = 1
$ws.Range("AL35").Value = 1
$ws.Range("AL36").Value = 1
$ws.Range("AL47").Value = 1
$ws.Range("AL48").Value = 1
$ws.Range("AL54").Value = 1
$ws.Range("AL61").Value = 1
$ws.Range("AL68").Value = 1
$ws.Range("AL75").Value = 1
$ws.Range("AL83").Value = 1
